$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (Joe Biden / Office of the Clerk) need to swap places so
# that "Office of the Clerk..." becomes row 2 and "Joe Biden..." becomes row 3.
# Sort the A2:E3 block by column A (ascending puts "Office..." before "Joe...").
$ws.Range("A1:E3").Sort($ws.Range("A2:A3"), 2)

# Range.Sort moves the cell values (and number/cell formats) with the rows,
# but the hyperlink objects stay anchored to their original cells, so their
# targets now point at the wrong row. Re-point them to match the text/URL
# that now lives in each row.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://web.archive.org/web/20120402071556/http://clerk.house.gov/floorsummary/floor.html?day=20090108")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.wsj.com/articles/joe-biden-decides-not-to-enter-presidential-race-1445444657")
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
